$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header columns: "_old" -> "_FV2404", "_new" -> "_FV2410"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace "_old$", "_FV2404"
        $newVal = $newVal -replace "_new$", "_FV2410"
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}

# Turn the header range into a proper Excel Table (ListObject)
$range = $ws.Range("A1:U77")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# Freeze the header row (split pane)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
